$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M")

# --- Row 8: period headers (shift one period left, append newest period) ---
$headers = @(
    "9 ماهه منتهی به 1399/09",
    "12 ماهه منتهی به 1399/12",
    "3 ماهه منتهی به 1400/03",
    "6 ماهه منتهی به 1400/06",
    "9 ماهه منتهی به 1400/09",
    "12 ماهه منتهی به 1400/12",
    "3 ماهه منتهی به 1401/03",
    "6 ماهه منتهی به 1401/06",
    "9 ماهه منتهی به 1401/09",
    "12 ماهه منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $headers[$i]
}

# --- Row 9: publish dates (shift one left, newest report date replaces, append newest) ---
$dates = @(
    "1400-10-29 (2)",
    "1401-02-07 (9)",
    "1401-04-26 (4)",
    "1401-08-28 (4)",
    "1401-10-29 (2)",
    "1402-02-13 (9)",
    "1401-04-26 (2)",
    "1401-08-28 (2)",
    "1401-10-29",
    "1402-02-13 (2)"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $dates[$i]
}

# --- Data rows 11-27: every row shifts one period left, newest period's data appended in M ---
$rowData = @{
    11 = @(2043796,2703782,783511,1657484,2753340,4106603,1120034,2406596,4092899,5968165)
    12 = @(-1239601,-1610046,-514765,-963756,-1484280,-2295951,-574049,-1261200,-2248671,-3398826)
    13 = @(804195,1093736,268746,693728,1269060,1810652,545985,1145396,1844228,2569339)
    14 = @(-202585,-228914,-31796,-64407,-102474,-158001,-39573,-93252,-166804,-256206)
    15 = @(0,0,0,0,0,0,0,0,0,0)
    16 = @(0,0,0,0,0,0,0,0,0,0)
    17 = @(601610,864822,236950,629321,1166586,1652651,506412,1052144,1677424,2313133)
    18 = @(0,0,-1625,-8149,-18756,-26174,-6727,-14242,-22576,-33600)
    19 = @(96792,118433,55538,116761,146657,182378,76660,103443,236593,255742)
    20 = @(698402,983255,290863,737933,1294487,1808855,576345,1141345,1891441,2535275)
    21 = @(-82429,-95721,-29107,-99653,-177637,-208189,-56023,-138997,-238607,-213046)
    22 = @(615973,887534,261756,638280,1116850,1600666,520322,1002348,1652834,2322229)
    23 = @(0,0,0,0,0,0,0,0,0,0)
    24 = @(615973,887534,261756,638280,1116850,1600666,520322,1002348,1652834,2322229)
    25 = @(4928,7100,2094,5106,8935,12805,2602,5012,8264,11611)
    26 = @(125000,125000,125000,125000,125000,125000,200000,200000,200000,200000)
    27 = @(3080,4438,1309,3191,5584,8003,2602,5012,8264,11611)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# --- Column widths: the "wide" (29-char) column shifts left by one period too ---
# ColumnWidth setter pads +5/6 internally, so subtract that to land on the exact
# target width written to the xlsx <col> element.
$pad = 5.0 / 6.0
$ws.Columns.Item(5).ColumnWidth  = 29 - $pad   # E -> 29
$ws.Columns.Item(6).ColumnWidth  = 28 - $pad   # F -> 28
$ws.Columns.Item(9).ColumnWidth  = 29 - $pad   # I -> 29
$ws.Columns.Item(10).ColumnWidth = 28 - $pad   # J -> 28
$ws.Columns.Item(13).ColumnWidth = 29 - $pad   # M -> 29
